$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPLKKPS001")

# Update USERID column (G) from "Putri" to numeric 32382
$ws.Range("G2").Value = 32382
$ws.Range("G3").Value = 32382
$ws.Range("G4").Value = 32382
$ws.Range("G5").Value = 32382

# Update PASSWORD column (H) from "bni1234/" to "bni1234"
$ws.Range("H2").Value = "bni1234"
$ws.Range("H3").Value = "bni1234"
$ws.Range("H4").Value = "bni1234"
$ws.Range("H5").Value = "bni1234"

# Update PREPARATION column (F) text to reflect new username/password
$ws.Range("F2").Value = "Username : 32382;`nPassword : bni1234;`nKode Bidang Usaha : 9;`nNama Bidang Usaha : Pertambangan"
$ws.Range("F3").Value = "Username : 32382;`nPassword : bni1234;`nKode Bidang Usaha : 9"
$ws.Range("F4").Value = "Username : 32382;`nPassword : bni1234;`nKode Bidang Usaha : 9;`nNama Bidang Usaha : Pariwisata"
$ws.Range("F5").Value = "Username : 32382;`nPassword : bni1234;`nKode Bidang Usaha : 9"

$ws.Range("F2").Select()
